$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new Price (column D) values below are plain digit strings like
# "100.44" that Excel would otherwise auto-convert to a Number (silently
# dropping significant trailing zeros, e.g. "301.60" -> 301.6). The source
# data is text (note other rows such as "46.643.76" have multiple dots and
# are already unambiguous text), so force a Text number format on exactly
# those cells before writing, keeping every other cell style untouched.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "46.643.76"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "2.273.98"
$ws.Range("E3").Value = "  -2.53%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "301.60"
$ws.Range("E5").Value = "  -1.81%  "
$ws.Range("D6").Value = "100.44"
$ws.Range("E6").Value = "  +2.05%  "
$ws.Range("D7").Value = "0.567"
$ws.Range("E7").Value = "  -1.60%  "
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").Value = "0.509"
$ws.Range("E9").Value = "  -5.44%  "
$ws.Range("D10").Value = "35.25"
$ws.Range("E10").Value = "  -2.75%  "
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("D12").Value = "7.09"
$ws.Range("E12").Value = "  -5.09%  "
$ws.Range("E13").Value = "  -1.68%  "
$ws.Range("D14").Value = "2.619.43"
$ws.Range("E14").Value = "  -2.56%  "
$ws.Range("D15").Value = "2.274.78"
$ws.Range("E15").Value = "  -2.54%  "
$ws.Range("D16").Value = "13.68"
$ws.Range("E16").Value = "  -3.35%  "
$ws.Range("D17").Value = "0.801"
$ws.Range("E17").Value = "  -3.96%  "
$ws.Range("D18").Value = "46.649.52"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").Value = "12.68"
$ws.Range("E19").Value = "  -2.19%  "
$ws.Range("E20").Value = "  +2.05%  "
$ws.Range("D21").Value = "5.87"
$ws.Range("E21").Value = "  -5.40%  "
$ws.Range("D22").Value = "66.00"
$ws.Range("E22").Value = "  -1.20%  "
$ws.Range("D23").Value = "249.63"
$ws.Range("E23").Value = "  +2.11%  "
$ws.Range("D24").Value = "2.81"
$ws.Range("E24").Value = "  -5.42%  "
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").Value = "  -5.58%  "
$ws.Range("D27").Value = "41.54"
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("D28").Value = "2.26"
$ws.Range("E28").Value = "  -1.14%  "
$ws.Range("D29").Value = "9.65"
$ws.Range("E29").Value = "  -1.95%  "
$ws.Range("D30").Value = "20.23"
$ws.Range("E30").Value = "  +0.38%  "
$ws.Range("E31").Value = "  +7.00%  "
$ws.Range("D32").Value = "3.40"
$ws.Range("E32").Value = "  +12.82%  "
$ws.Range("D33").Value = "147.21"
$ws.Range("E33").Value = "  -3.40%  "
$ws.Range("D34").Value = "5.40"
$ws.Range("E34").Value = "  -5.62%  "
$ws.Range("E35").Value = "  -4.31%  "
$ws.Range("E36").Value = "  +6.19%  "
$ws.Range("E37").Value = "  -2.65%  "
$ws.Range("D38").Value = "15.83"
$ws.Range("E38").Value = "  +12.71%  "
$ws.Range("D39").Value = "1.70"
$ws.Range("E39").Value = "  -6.23%  "
$ws.Range("E40").Value = "  -3.99%  "
$ws.Range("D41").Value = "0.0297"
$ws.Range("E41").Value = "  -6.72%  "
$ws.Range("E42").Value = "  -7.56%  "
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("D44").Value = "93.26"
$ws.Range("E44").Value = "  +15.36%  "
$ws.Range("D45").Value = "1.793.98"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("E46").Value = "  -4.72%  "
$ws.Range("D47").Value = "71.28"
$ws.Range("E47").Value = "  -4.46%  "
$ws.Range("E48").Value = "  -6.54%  "
$ws.Range("D49").Value = "4.81"
$ws.Range("E49").Value = "  -1.76%  "
$ws.Range("D50").Value = "95.15"
$ws.Range("E50").Value = "  -3.07%  "
$ws.Range("D51").Value = "7.92"
$ws.Range("E51").Value = "  -1.03%  "
